$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.124.38"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "3.742.46"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "601.27"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "167.29"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "3.740.70"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("E10").Value = "  +3.69%  "

$ws.Range("E11").Value = "  +0.58%  "

$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("D15").Value = "4.368.30"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "3.742.99"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").Value = "69.072.81"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "7.35"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("E20").Value = "  -1.52%  "

$ws.Range("E21").Value = "  +10.99%  "

$ws.Range("D22").Value = "492.37"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("E24").Value = "  +8.49%  "

$ws.Range("D25").Value = "85.03"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  -0.68%  "

$ws.Range("D27").Value = "12.24"
$ws.Range("E27").Value = "  -0.90%  "

$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +1.18%  "

$ws.Range("D31").Value = "8.16"
$ws.Range("E31").Value = "  +2.58%  "

$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("D34").Value = "3.888.24"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("E39").Value = "  +2.17%  "

$ws.Range("E40").Value = "  +4.74%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").Value = "3.04"
$ws.Range("E42").Value = "  +6.18%  "

$ws.Range("D43").Value = "48.79"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("E44").Value = "  +0.53%  "

$ws.Range("D45").Value = "423.37"
$ws.Range("E45").Value = "  -2.64%  "

$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("D48").Value = "40.03"
$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("D49").Value = "141.90"
$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("D50").Value = "2.782.41"
$ws.Range("E50").Value = "  +1.50%  "

$ws.Range("E51").Value = "  -0.05%  "

# Row 35 / 36: Hedera and RenzoRestakedETH swap ranking positions with updated values
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.675.90"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").Value = "  -0.23%  "

